$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel COM colour values are packed as 0x00BBGGRR (BGR), matching these OOXML RGB colours:
$colGreen = 4113922   # FF02C63E -> story card already completed
$colGray  = 10066329  # FF999999 -> story card not yet started
$colBlue  = 16755456  # FF00ABFF -> story card in progress
$xlPasteFormats = -4122

# Apply the look of an existing "story card" cell (grey fill, thin grey border,
# thick coloured left edge, Arial font, wrap + top-align + indent) to $range, then
# recolour its left edge to signal status.
function Set-CardFormat($range, $color) {
    $ws.Range("B9").Copy() | Out-Null
    $range.PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Application.CutCopyMode = 0
    if ($color -eq "gray") { $range.Borders.Item(7).Color = $colGray }
    elseif ($color -eq "blue") { $range.Borders.Item(7).Color = $colBlue }
    else { $range.Borders.Item(7).Color = $colGreen }
}

# Apply the look of an existing "Release/epic" divider row (dashed bottom border) to $range.
# These rows use the sheets default row height (no explicit <row ht=.../> override), so
# AutoFit() resets the row back to standard height + drops any stale custom-height flag
# left over from content that previously lived at that row number.
function Set-HeaderFormat($range, $rowNum) {
    $ws.Range("B7:H7").Copy() | Out-Null
    $range.PasteSpecial($xlPasteFormats) | Out-Null
    $ws.Application.CutCopyMode = 0
    $ws.Rows.Item($rowNum).AutoFit()
}

# Rows 21-45 held the old "Release 1/2/3 MVP" backlog; it is being reflowed (new rows
# inserted, stories re-ordered) and extended through row 53 with a new "Release 4" lane.
# Clearing first guarantees no stale values/styles are left behind at addresses that
# do not get rewritten below (e.g. old divider row 41 spanned B:H, new row 41 is a
# single story cell).
$ws.Range("B21:H53").Clear()

# --- Row 21 ---
$ws.Range("B21").Value = "Hifi prototype of web app"
$ws.Range("D21").Value = "Update Problem Definition"
$ws.Range("F21").Value = "Outline CICD pipeline"
Set-CardFormat $ws.Range("B21") "green"
Set-CardFormat $ws.Range("D21") "green"
Set-CardFormat $ws.Range("F21") "green"
$ws.Rows.Item(21).RowHeight = 72

$ws.Rows.Item(22).RowHeight = 10

# --- Row 23 ---
$ws.Range("B23").Value = "Continue Architecture Research"
$ws.Range("D23").Value = "Define coding standards and guidelines"
Set-CardFormat $ws.Range("B23") "green"
Set-CardFormat $ws.Range("D23") "green"
$ws.Rows.Item(23).RowHeight = 72

$ws.Rows.Item(24).RowHeight = 10

# --- Row 25 ---
$ws.Range("B25").Value = "Release 1 MVP"
Set-HeaderFormat $ws.Range("B25:H25") 25

$ws.Rows.Item(26).RowHeight = 10

# --- Row 27 ---
$ws.Range("B27").Value = "Update USM"
$ws.Range("D27").Value = "Research available libraries"
$ws.Range("F27").Value = "Add Reporting Page"
$ws.Range("H27").Value = "Add Additional Report"
Set-CardFormat $ws.Range("B27") "green"
Set-CardFormat $ws.Range("D27") "green"
Set-CardFormat $ws.Range("F27") "green"
Set-CardFormat $ws.Range("H27") "green"
$ws.Rows.Item(27).RowHeight = 72

$ws.Rows.Item(28).RowHeight = 10

# --- Row 29 ---
$ws.Range("B29").Value = "Setup CICD for front end"
$ws.Range("D29").Value = "Add Recycling Overview Page"
Set-CardFormat $ws.Range("B29") "green"
Set-CardFormat $ws.Range("D29") "green"
$ws.Rows.Item(29).RowHeight = 72

$ws.Rows.Item(30).RowHeight = 10

# --- Row 31 ---
$ws.Range("B31").Value = "Create Web App Project"
Set-CardFormat $ws.Range("B31") "green"
$ws.Rows.Item(31).RowHeight = 72

$ws.Rows.Item(32).RowHeight = 10

# --- Row 33 ---
$ws.Range("B33").Value = "Install necessary libraries"
Set-CardFormat $ws.Range("B33") "green"
$ws.Rows.Item(33).RowHeight = 72

$ws.Rows.Item(34).RowHeight = 10

# --- Row 35 ---
$ws.Range("B35").Value = "Release 2 MVP"
Set-HeaderFormat $ws.Range("B35:H35") 35

$ws.Rows.Item(36).RowHeight = 10

# --- Row 37 ---
$ws.Range("B37").Value = "Add Login Page"
$ws.Range("D37").Value = "Add heat map of contaminant disposal"
$ws.Range("F37").Value = "Add Aggregate Recycling Data"
$ws.Range("H37").Value = "Add Problem Neighbourhood information"
Set-CardFormat $ws.Range("B37") "gray"
Set-CardFormat $ws.Range("D37") "gray"
Set-CardFormat $ws.Range("F37") "green"
Set-CardFormat $ws.Range("H37") "green"
$ws.Rows.Item(37).RowHeight = 72

$ws.Rows.Item(38).RowHeight = 10

# --- Row 39 ---
$ws.Range("B39").Value = "Create API Endpoints"
Set-CardFormat $ws.Range("B39") "blue"
$ws.Rows.Item(39).RowHeight = 72

$ws.Rows.Item(40).RowHeight = 10

# --- Row 41 ---
$ws.Range("B41").Value = "Tie in API calls and database storage"
Set-CardFormat $ws.Range("B41") "gray"
$ws.Rows.Item(41).RowHeight = 72

$ws.Rows.Item(42).RowHeight = 10

# --- Row 43 ---
$ws.Range("B43").Value = "Update CICD pipeline"
Set-CardFormat $ws.Range("B43") "gray"
$ws.Rows.Item(43).RowHeight = 72

$ws.Rows.Item(44).RowHeight = 10

# --- Row 45 ---
$ws.Range("B45").Value = "Release 3 MVP"
Set-HeaderFormat $ws.Range("B45:H45") 45

$ws.Rows.Item(46).RowHeight = 10

# --- Row 47 ---
$ws.Range("B47").Value = "Add additional API endpoints"
$ws.Range("D47").Value = "Add Zoom Function"
$ws.Range("F47").Value = "Add Export Function"
$ws.Range("H47").Value = "Add Column Hiding"
Set-CardFormat $ws.Range("B47") "gray"
Set-CardFormat $ws.Range("D47") "gray"
Set-CardFormat $ws.Range("F47") "gray"
Set-CardFormat $ws.Range("H47") "gray"
$ws.Rows.Item(47).RowHeight = 72

$ws.Rows.Item(48).RowHeight = 10

# --- Row 49 ---
$ws.Range("B49").Value = "Add Login Functionality"
$ws.Range("D49").Value = "Add tool tip on hover"
Set-CardFormat $ws.Range("B49") "gray"
Set-CardFormat $ws.Range("D49") "gray"
$ws.Rows.Item(49).RowHeight = 72

$ws.Rows.Item(50).RowHeight = 10

# --- Row 51 ---
$ws.Range("B51").Value = "Release 4 Future Release MVP"
Set-HeaderFormat $ws.Range("B51:H51") 51

$ws.Rows.Item(52).RowHeight = 10

# --- Row 53 ---
$ws.Range("B53").Value = "Add Test and Staging server to CICD pipeline"
Set-CardFormat $ws.Range("B53") "gray"
$ws.Rows.Item(53).RowHeight = 72

# Sheet now spans through column H, row 53 (was row 45).
Write-Host "Edit complete"